$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Type" column header in G1, matching the style of the other header cells
$ws.Range("G1").Value = "Type"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = "Type"

# Fill rows 2-7 in column G with "P" (parish)
$ws.Range("G2:G7").Value = "P"

# Update the selection to match the diff (G8)
$ws.Range("G8").Select()
